# Updated symbol list on Sat Dec 31 18:37:45 UTC 2022 with GitHub Actions
#
# Cryptocurrency prices/volumes are stored as plain text (e.g. "246.56",
# "0.63%") in columns D (Price) and E (Volume 1h). Assigning a numeric- or
# percentage-looking string directly would make Excel auto-convert the cell
# to a Number, so each such cell is pre-formatted as Text ("@") before the
# new value is written, keeping it a text value exactly like the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "246.62"
Set-TextValue "E2" "0.54%"

# Row 3 - OKB
Set-TextValue "D3" "26.28"
Set-TextValue "E3" "5.33%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.084"
Set-TextValue "E4" "1.10%"

# Row 5 - Cronos
Set-TextValue "D5" "0.05599"
Set-TextValue "E5" "-0.38%"

# Row 6 - KuCoinToken
Set-TextValue "D6" "6.486"
Set-TextValue "E6" "-0.78%"

# Row 7 - MXToken
Set-TextValue "D7" "0.8132"
Set-TextValue "E7" "0.44%"

# Row 8 - FTXToken
Set-TextValue "D8" "0.8472"
Set-TextValue "E8" "1.00%"

# Row 9 - BitrueCoin
Set-TextValue "D9" "0.02846"
Set-TextValue "E9" "0.45%"

# Row 10 - BitMartToken
Set-TextValue "D10" "0.09394"
Set-TextValue "E10" "-0.15%"

# Row 11 - BitForexToken
Set-TextValue "D11" "0.001513"
Set-TextValue "E11" "-0.79%"

# Row 12 - One
Set-TextValue "D12" "0.0005992"
Set-TextValue "E12" "0.74%"

# Row 13 - TigerCash
Set-TextValue "D13" "0.006197"
Set-TextValue "E13" "-0.89%"

# Row 14 - LEO
Set-TextValue "D14" "3.607"
Set-TextValue "E14" "3.10%"

# Row 15 - GateToken
Set-TextValue "D15" "3.013"
Set-TextValue "E15" "0.83%"

# Row 16 - BTSEToken
Set-TextValue "E16" "0.05%"

# Row 18 - WazirX
Set-TextValue "D18" "0.1340"
Set-TextValue "E18" "0.21%"

# Row 19 - MandalaExchangeToken
Set-TextValue "D19" "0.07007"
Set-TextValue "E19" "0.96%"

# Row 20 - LiechtensteinCryptoassetsExchange
Set-TextValue "D20" "0.03176"
Set-TextValue "E20" "-3.45%"

# Row 21 - ProBitToken
Set-TextValue "E21" "0.46%"

# Row 22 - MCDex
Set-TextValue "D22" "3.745"
Set-TextValue "E22" "-0.03%"

# Row 23 - CoinExToken
Set-TextValue "D23" "0.04657"
Set-TextValue "E23" "-0.58%"

# Row 24 - ZBToken
Set-TextValue "E24" "-1.41%"

# Row 25 - BitKan
Set-TextValue "D25" "0.001245"
Set-TextValue "E25" "0.11%"

# Row 26 - HotbitToken
Set-TextValue "E26" "1.44%"

# Row 27 - NitroEx
Set-TextValue "D27" "0.00009602"
Set-TextValue "E27" "-0.95%"

# Row 28 - UpBots
Set-TextValue "D28" "0.0001939"
Set-TextValue "E28" "-0.04%"

# Row 40 - IDEX
Set-TextValue "D40" "0.03670"
Set-TextValue "E40" "1.15%"

# Row 41 - now BKEXToken (was KickToken)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1363"
Set-TextValue "E41" "29.59%"

# Row 42 - now CEJI (was BKEXToken)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002661"
Set-TextValue "E42" "5.99%"

# Row 43 - now KickToken (was CEJI)
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006231"
Set-TextValue "E43" "-0.55%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008759"
Set-TextValue "E44" "4.77%"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005298"
Set-TextValue "E45" "0.44%"

# Row 46 - Kangarootoken
Set-TextValue "E46" "0.10%"

# Row 47 - CoinbaseStockToken
Set-TextValue "E47" "-39.96%"

# Row 48 - BOLO
Set-TextValue "D48" "0.002484"
Set-TextValue "E48" "21.25%"

# Row 49 - CryptobidCoin
Set-TextValue "D49" "0.00002101"
Set-TextValue "E49" "0.10%"

# Row 50 - SpecialPowerGold
Set-TextValue "D50" "0.0002001"
Set-TextValue "E50" "0.10%"
